# feat: evolução tela login com fundo e imagens, correção erro relacional
#
# Updates the relatorios sheet: two report rows get renamed/repointed to
# new PowerBI dashboards, and the two remaining PowerBI links are swapped
# out for the (placeholder) Google URLs used by the new login screen
# background/image flow. Hyperlink targets are refreshed to match, column B
# is widened to fit the new longer report names, and the active selection
# is left on E3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh the link column first (C2:C5) --------------------------------
# Row 2 / Row 3 point at the two new PowerBI dashboards; Row 4 / Row 5 keep
# their report names but their links are corrected (relational error fix)
# to point at the new targets.
$ws.Range("C2").Value = "https://app.powerbi.com/view?r=eyJrIjoiMTU3YWQxYjktYTI5MC00OTFmLWFlYzItYmZlMGZiZTRjNmVjIiwidCI6ImM1M2UwMWZmLTQ5MjItNDAzYy1iZWE1LTQzOTZjMmUxMmQ5NyJ9"
$ws.Range("C3").Value = "https://app.powerbi.com/view?r=eyJrIjoiZTg4MTQwNDEtNWFkYS00ZjY3LWEyMzItMjNlODczZGEyMDE3IiwidCI6ImM1M2UwMWZmLTQ5MjItNDAzYy1iZWE1LTQzOTZjMmUxMmQ5NyJ9"
$ws.Range("C4").Value = "https://www.google.com/?hl=pt_BR"
$ws.Range("C5").Value = "https://www.google.com/imghp?hl=pt-BR&ogbl"

# --- Row 2: "Painel Recuperação" -> "Painel CIG Rentabilidade" ------------
$ws.Range("B2").Value = "Painel CIG Rentabilidade"

# --- Row 3: "Fluxo de Caixa" -> "Painel CIG Saldo Bancário" ---------------
$ws.Range("B3").Value = "Painel CIG Saldo Bancário"

# --- Refresh the hyperlink objects so C2:C5 point at the new addresses ---
# (the worksheet Hyperlinks collection has to be rebuilt - editing
# .Address in place on an existing Hyperlink object does not repoint it).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C2"), $ws.Range("C2").Value())
$ws.Hyperlinks.Add($ws.Range("C3"), $ws.Range("C3").Value())
$ws.Hyperlinks.Add($ws.Range("C4"), $ws.Range("C4").Value())
$ws.Hyperlinks.Add($ws.Range("C5"), $ws.Range("C5").Value())

# Re-adding hyperlinks reapplies direct formatting; snap C2:C5 back onto
# the shared "Hiperlink" cell style instead of leaving ad-hoc formatting.
$ws.Range("C2:C5").Style = "Hiperlink"

# --- Column B needs to widen to fit the longer report names --------------
$ws.Columns("B:B").ColumnWidth = 22.43

# --- Leave the selection on E3, matching the saved workbook state --------
$ws.Range("E3").Select()
